$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D28: 568 -> 492
$ws.Range("D28").Value = 492

# Row 37 currently holds: 2970482 / Pril Isis Ultra Power 650ml / 12 / 168.
# That whole record moves down, unchanged, into a brand-new row 39. Populate
# row 39 BEFORE touching row 37 itself, so the "2970482" shared string stays
# referenced the entire time (otherwise it gets dropped once row 37 no
# longer points at it, and re-inserting it later collides/renumbers).
$ws.Range("B39").Value = "Pril Isis Ultra Power 650ml"
$ws.Range("C39").Value = 12
$ws.Range("D39").Value = 168

# A39 = "2970482" must remain text (it looks like a plain number, and a
# bare .Value assignment would silently coerce it to numeric). Build it as
# a text formula first, then flatten formula -> literal value via
# copy/paste-values - this keeps it a genuine shared-string cell without
# Excel's numeric coercion and without minting a new, unused number-format
# style (which a NumberFormat="@" / leading-apostrophe trick would do).
$ws.Range("A39").Formula = "=""2970482"""
$ws.Range("A39").Copy()
$ws.Range("A39").PasteSpecial(-4163)   # xlPasteValues

# Copy row 38's formatting (the thin-border cell style) down onto row 39.
$ws.Range("A38:D38").Copy()
$ws.Range("A39:D39").PasteSpecial(-4122)   # xlPasteFormats

# Row 37 itself is replaced with a new item: 2952095 / LE CHAT power gel 4L / 3 / 1190
# A37 also loses its cell formatting (no more border style) in this edit.
$ws.Range("A37").Value = 2952095
$ws.Range("A37").ClearFormats()
$ws.Range("B37").Value = "LE CHAT power gel 4L"
$ws.Range("C37").Value = 3
$ws.Range("D37").Value = 1190

# Reflect the newly added row in the sheet's selection/scroll state.
$null = $ws.Range("A26").Select()
$null = $ws.Range("D29").Select()
